$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing row (21) currently has the "final row" date-only format.
# That format now moves to the new last row (22), so row 21's date cell
# reverts to the regular date-time format used by the rest of the column.
$ws.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 22.
$ws.Range("A22").Value = 45606
$ws.Range("B22").Value = 52
$ws.Range("C22").Value = 46
$ws.Range("D22").Value = 53

# The new last row gets the distinctive date-only format that row 21 used
# to have.
$ws.Range("A22").NumberFormat = "YYYY-MM-DD"
